$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# This workbook originally had quarterly financial data in columns D:K
# (most-recent-quarter first). The update adds two newer quarters of
# data, which get inserted as new columns D and E, pushing the existing
# D:K data to F:M.
# -----------------------------------------------------------------------

# 1. Insert two new blank columns at D:E - this shifts old D:K to F:M
#    and keeps all existing formatting/styles on the shifted cells intact.
$ws.Columns("D:E").Insert()

# 2. Stamp column D:E with the same per-row number formats as column F
#    (the former column D) so date rows keep the date format and data
#    rows keep the numeric format, without minting new style entries.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Populate the new D/E columns with the newly reported quarter values.
#    Each triple is (row, newD, newE).
$newData = @(
    @(7, 43465, 43373),
    @(8, 254400, 244000),
    @(9, "NA", "NA"),
    @(10, "NA", "NA"),
    @(11, $null, $null),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 0, 0),
    @(15, 0, 0),
    @(16, $null, $null),
    @(17, 88600, 79100),
    @(18, 165800, 164900),
    @(19, $null, $null),
    @(20, -123000, -91500),
    @(21, 46700, 78500),
    @(22, 0, 0),
    @(23, 42800, 73400),
    @(24, 11700, 19200),
    @(25, 0, 0),
    @(26, 31000, 54200),
    @(27, 31000, 54200),
    @(28, 0, 0),
    @(29, 2300, "NA"),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 123000, 91500),
    @(33, 33300, 54200),
    @(34, 0, 0),
    @(35, 33300, 54200),
    @(38, 43465, 43373),
    @(39, $null, $null),
    @(40, $null, $null),
    @(41, 196900, 210600),
    @(42, 260200, 242400),
    @(43, 0, 0),
    @(44, 0, 0),
    @(45, 0, 0),
    @(46, 0, 0),
    @(47, 0, 0),
    @(48, 177400, 175400),
    @(49, 99100, 99800),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 104400, 135500),
    @(53, 0, 0),
    @(54, 26229000, 25518500),
    @(55, $null, $null),
    @(56, $null, $null),
    @(57, 0, 0),
    @(58, 0, 0),
    @(59, 0, 0),
    @(60, 0, 0),
    @(61, 5435700, 4853800),
    @(62, 0, 0),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 23223700, 22483300),
    @(67, $null, $null),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 1173900, 1172600),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 3005300, 3035200),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 33300, 54200),
    @(82, $null, $null),
    @(83, 4000, 5100),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 42600, 80000),
    @(90, $null, $null),
    @(91, -5500, -3400),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -720800, -154500),
    @(95, $null, $null),
    @(96, -32000, -26700),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, 664400, 89200),
    @(101, 0, 0),
    @(102, -13700, 14600)

)

foreach ($entry in $newData) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
}

Write-Output "Applied ISBC quarterly financials update"
